$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jaana")

# Copy the date-cell formatting from the last filled row (row 18) down onto the
# two new rows so the new date cells pick up the same date number-format /
# border style the sheet already uses for the other entries, instead of
# Excel inventing a brand new style.
$ws.Range("A18").Copy()
$ws.Range("A19:A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New timesheet entries for Jaana - "Koodikorjauksia Customer/Invoice" work.
$ws.Range("A19").Value = 44994
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "Koodikorjauksia Customer/Invoice"

$ws.Range("A20").Value = 44995
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = "Dokumentaation päivitys, useamman asiakkaan laskujen näyttäminen listassa."

# Make "Jaana" the active/selected sheet with A21 selected (matches the
# author re-saving the workbook with focus back on her own tab).
$ws.Activate() | Out-Null
$ws.Range("A21").Select() | Out-Null
